# ARKCORR-18: Removing assignee when object moves through queues
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Functions cell (D16): addOrUpdateParticipant now resolves "id" via a SpEL
#     expression instead of taking it verbatim ---
$ws.Range("D16").Value = 'function String dateFormat(String fmt)
{
  return LocalDate.now().format(DateTimeFormatter.ofPattern(fmt));
}

function Date toDate(LocalDate localDate)
{
    return Date.from(localDate.atStartOfDay(ZoneId.systemDefault()).toInstant());
}

function Boolean evalSpring(String expression, Object obj)
{
    ExpressionParser ep = new SpelExpressionParser();
    Expression exp = ep.parseExpression(expression);
    EvaluationContext ec = new StandardEvaluationContext();

    Boolean evaluated = exp.getValue(ec, obj, Boolean.class);
 
    return evaluated;
}
function void addOrUpdateParticipant(AcmAssignedObject obj, String type, String expression) 
{
    ExpressionParser ep = new SpelExpressionParser();
    Expression exp = ep.parseExpression(expression);
    EvaluationContext ec = new StandardEvaluationContext();
    String id = exp.getValue(ec, obj, String.class);
    if  (obj.getParticipants() != null && type != null)
    {
        boolean updated = false;
        for (AcmParticipant p :  obj.getParticipants())
        {
             if (type.equals(p.getParticipantType()))
             {
                  p.setParticipantLdapId(id);
                  updated = true;
                  break;
             }
        }
        if (!updated)
        {
            AcmParticipant p = new AcmParticipant();
            p.setParticipantType(type);
            p.setParticipantLdapId(id);
 
            obj.getParticipants().add(p);
        }
    }
}'

# --- Owning-group actions (E30:E34): wrap the literal group name in
#     new String(...) so each rule call gets a fresh String instance ---
$ws.Range("E30").Value = 'owning group, new String(''ExecSec Intake'')'
$ws.Range("E31").Value = 'owning group, new String(''Analyst - AG'')'
$ws.Range("E32").Value = 'owning group, new String(''Supervisor - AG'')'
$ws.Range("E33").Value = 'owning group, new String(''ExecSec Approval'')'
$ws.Range("E34").Value = 'owning group, new String(''ExecSec Release'')'

# --- New "Set Queue Enter Date *" rows (35:39): clear the assignee whenever
#     the case file enters a new queue ---
$ws.Range("E35").Value = 'assignee, new String('''')'
$ws.Range("E36").Value = 'assignee, new String('''')'
$ws.Range("E37").Value = 'assignee, new String('''')'
$ws.Range("E38").Value = 'assignee, new String('''')'
$ws.Range("E39").Value = 'assignee, new String('''')'

# --- Restore view state: active cell / scroll focus moved to D16 ---
$ws.Range("D16").Select()

# --- Touch page setup so orientation stays explicit on the sheet (portrait) ---
$ws.PageSetup.Orientation = 1
